$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheet name: $($ws.Name)"
Write-Host "UsedRange: $($ws.UsedRange.Address())"
